$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update yoy_schools, yoy_authorities, yoy_users for rows 14-20 (2025-07 stats update)
# Values were multiplied by 100 (converted from fraction to percentage number)
$updates = @{
    "F14" = 2.940153096729303
    "G14" = 5.961754780652417
    "H14" = 23.28061250163025

    "F15" = 3.658536585365857
    "G15" = 6.764374295377684
    "H15" = 24.85887932178075

    "F16" = 4.054289194362282
    "G16" = 6.877113866967299
    "H16" = 25.07756835683654

    "F17" = 5.86376404494382
    "G17" = 3.205128205128216
    "H17" = 18.87096770378025

    "F18" = 6.092436974789917
    "G18" = 3.311965811965822
    "H18" = 18.41667687390272

    "F19" = 6.339254615116685
    "G19" = 3.201707577374591
    "H19" = 24.62859203576528

    "F20" = 6.184142338918641
    "G20" = 3.503184713375807
    "H20" = 25.58277891171774
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
